$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.572.46"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").Value = "1.979.29"
$ws.Range("E3").Value = "  -3.74%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.58%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.359"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0733"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "

$ws.Range("E12").Value = "  -2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.941"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "2.264.46"
$ws.Range("E15").Value = "  -4.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("D17").Value = "1.979.90"
$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.77%  "

$ws.Range("D19").Value = "35.544.74"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.42%  "

$ws.Range("E26").Value = "  -2.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0961"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0590"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.86%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  -4.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.72%  "

$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("E42").Value = "  -2.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "91.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0883"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "

$ws.Range("D48").Value = "1.371.85"
$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.54%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
